$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7581236958503723
$ws.Range("B1").Value = 1.422454953193665
$ws.Range("C1").Value = 5.391351699829102
$ws.Range("D1").Value = 3.174206733703613
$ws.Range("E1").Value = 1.514097332954407
